$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.064883708953857
$ws.Range("B1").Value = 2.654764175415039
$ws.Range("C1").Value = 8.938697814941406
$ws.Range("D1").Value = 2.052101612091064
$ws.Range("E1").Value = 1.162778973579407
